$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 358.58334
$ws.Range("I28").Value = 358.58334
$ws.Range("K28").Value = 358.58334
$ws.Range("M28").Value = 126.41666
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 100
$ws.Range("K29").Value = 300
$ws.Range("M29").Value = -19
$ws.Range("H32").Value = 1396.55
$ws.Range("J32").Value = 1686.4166
$ws.Range("L32").Value = 1686.4166
$ws.Range("N32").Value = -2338.4166
$ws.Range("H33").Value = 269.25
$ws.Range("I33").Value = 192.33333
$ws.Range("K33").Value = 192.33333
$ws.Range("M33").Value = 36.66667000000001
$ws.Range("H40").Value = 3122
$ws.Range("I40").Value = 4024.625
$ws.Range("J40").Value = 2399.9
$ws.Range("K40").Value = 4024.625
$ws.Range("L40").Value = 2399.9
$ws.Range("M40").Value = -3849.625
$ws.Range("N40").Value = -2749.9
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H105").Value = 888.5
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H116").Value = 1966.6666
$ws.Range("I116").Value = 1966.6666
$ws.Range("K116").Value = 1966.6666
$ws.Range("M116").Value = 1475.3334
$ws.Range("H132").Value = 1320
$ws.Range("I132").Value = 1320
$ws.Range("K132").Value = 3960
$ws.Range("M132").Value = -1430
$ws.Range("H138").Value = 3538.4517
$ws.Range("J138").Value = 6861
$ws.Range("L138").Value = 20583
$ws.Range("N138").Value = -30863

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 674.3143
$ws.Range("I74").Value = 691.15625
$ws.Range("K74").Value = 691.15625
$ws.Range("M74").Value = 182.84375
$ws.Range("H77").Value = 674.3143
$ws.Range("I77").Value = 691.15625
$ws.Range("K77").Value = 3455.78125
$ws.Range("M77").Value = 912.21875
$ws.Range("H97").Value = 775.17645
$ws.Range("I97").Value = 773.6
$ws.Range("K97").Value = 773.6
$ws.Range("M97").Value = -277.6
$ws.Range("H132").Value = 35116.5
$ws.Range("I132").Value = 2733
$ws.Range("K132").Value = 8199
$ws.Range("M132").Value = -5669

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4090.5557
$ws.Range("I86").Value = 4188.857
$ws.Range("J86").Value = 3746.5
$ws.Range("K86").Value = 4188.857
$ws.Range("L86").Value = 3746.5
$ws.Range("M86").Value = -3065.857
$ws.Range("N86").Value = -5992.5
$ws.Range("H89").Value = 4090.5557
$ws.Range("I89").Value = 4188.857
$ws.Range("J89").Value = 3746.5
$ws.Range("K89").Value = 20944.285
$ws.Range("L89").Value = 18732.5
$ws.Range("M89").Value = -15328.285
$ws.Range("N89").Value = -29964.5
$ws.Range("H107").Value = 2058.111
$ws.Range("I107").Value = 2305.5
$ws.Range("J107").Value = 79
$ws.Range("K107").Value = 2305.5
$ws.Range("L107").Value = 79
$ws.Range("M107").Value = -385.5
$ws.Range("N107").Value = -3919
$ws.Range("H134").Value = 1126.3667
$ws.Range("I134").Value = 872.15
$ws.Range("J134").Value = 1634.8
$ws.Range("K134").Value = 2616.45
$ws.Range("L134").Value = 4904.4
$ws.Range("M134").Value = -81.44999999999982
$ws.Range("N134").Value = -9974.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5002933
$ws.Range("I6").Value = 5002933
$ws.Range("K6").Value = 5002933
$ws.Range("M6").Value = -5002820
$ws.Range("H58").Value = 2544.8462
$ws.Range("I58").Value = 1316.5555
$ws.Range("K58").Value = 1316.5555
$ws.Range("M58").Value = -1113.5555
$ws.Range("H108").Value = 79000
$ws.Range("J108").Value = 79000
$ws.Range("L108").Value = 79000
$ws.Range("N108").Value = -86680
$ws.Range("H134").Value = 2972
$ws.Range("I134").Value = 2956.4285
$ws.Range("K134").Value = 8869.2855
$ws.Range("M134").Value = -6334.2855
$ws.Range("H136").Value = 2544.8462
$ws.Range("I136").Value = 1316.5555
$ws.Range("K136").Value = 3949.6665
$ws.Range("M136").Value = -1399.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 4239.933
$ws.Range("J62").Value = 3584.5386
$ws.Range("L62").Value = 10753.6158
$ws.Range("N62").Value = -12125.6158
$ws.Range("H65").Value = 4239.933
$ws.Range("J65").Value = 3584.5386
$ws.Range("L65").Value = 32260.8474
$ws.Range("N65").Value = -39124.8474
$ws.Range("H82").Value = 13868
$ws.Range("I82").Value = 9998.75
$ws.Range("J82").Value = 17737.25
$ws.Range("K82").Value = 29996.25
$ws.Range("L82").Value = 53211.75
$ws.Range("M82").Value = -29590.25
$ws.Range("N82").Value = -54023.75
$ws.Range("H85").Value = 13868
$ws.Range("I85").Value = 9998.75
$ws.Range("J85").Value = 17737.25
$ws.Range("K85").Value = 29996.25
$ws.Range("L85").Value = 53211.75
$ws.Range("M85").Value = -28592.25
$ws.Range("N85").Value = -56019.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4043.5
$ws.Range("I126").Value = 4202.2
$ws.Range("J126").Value = 3250
$ws.Range("K126").Value = 12606.6
$ws.Range("L126").Value = 9750
$ws.Range("M126").Value = -10136.6
$ws.Range("N126").Value = -14690

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 1258249.2
$ws.Range("J43").Value = 1674999.4
$ws.Range("L43").Value = 1674999.4
$ws.Range("N43").Value = -1675385.4
$ws.Range("H46").Value = 837
$ws.Range("I46").Value = 700.5
$ws.Range("J46").Value = 1110
$ws.Range("K46").Value = 700.5
$ws.Range("L46").Value = 1110
$ws.Range("M46").Value = -512.5
$ws.Range("N46").Value = -1486
$ws.Range("H104").Value = 34916.668
$ws.Range("J104").Value = 34916.668
$ws.Range("L104").Value = 34916.668
$ws.Range("N104").Value = -41904.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 9996.666999999999
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H37").Value = 25000
$ws.Range("I37").Value = 30000
$ws.Range("K37").Value = 30000
$ws.Range("M37").Value = -29797
$ws.Range("H105").Value = 42807.5
$ws.Range("J105").Value = 42807.5
$ws.Range("L105").Value = 42807.5
$ws.Range("N105").Value = -49795.5
$ws.Range("H132").Value = 4413.3
$ws.Range("J132").Value = 5865.8887
$ws.Range("L132").Value = 17597.6661
$ws.Range("N132").Value = -22657.6661
$ws.Range("H136").Value = 1341.1333
$ws.Range("I136").Value = 932.0769
$ws.Range("K136").Value = 2796.2307
$ws.Range("M136").Value = -246.2307000000001
